$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 52) to the tracking sheet, following the same
# layout as the existing rows: 日付 (date as text), 曜日 (weekday as text),
# 時刻 (hour, numeric), ランキング (ranking, numeric).

# Force column A to be treated as plain text so the date-like string
# "2025/10/02" is not auto-converted into a real date serial value.
$ws.Cells.Item(52, 1).NumberFormat = "@"
$ws.Cells.Item(52, 1).Value = "2025/10/02"
# Reset the style back to the default "Normal" style so no extra
# number-format/style index is left attached to the cell.
$ws.Cells.Item(52, 1).Style = "Normal"

$ws.Cells.Item(52, 2).Value = "木"
$ws.Cells.Item(52, 3).Value = 20
$ws.Cells.Item(52, 4).Value = 3
